# "added JAXB export to file"
#
# Adds a missing student (Сидоров Е. Г.) to the "Студенты" sheet, fixes the
# course-attendance value that was mistakenly stored as a 2% percentage
# instead of a plain count of 2, fills in the missing abbreviation for the
# Tambov medical university on the "Университеты" sheet, and corrects a
# typo'd founding year (2025 -> 2003) for the Voronezh university.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Студенты
$ws2 = $wb.Worksheets.Item(2)   # Университеты

# --- Студенты sheet -------------------------------------------------------

# New student row for course id "0003-high" (row 7), previously had no name.
$ws1.Range("B7").Value = "Сидоров Е. Г."

# C8 held 0.02 formatted as a percentage (displaying "2%"); it should simply
# be the integer course number 2, matching the formatting used by the other
# rows in that column.
$ws1.Range("C8").NumberFormat = "0.00"
$ws1.Range("C8").Value = 2

# --- Университеты sheet ----------------------------------------------------

# Missing abbreviation for "Тамбовский Университет Медицины".
$ws2.Range("C5").Value = "ТУМ"

# Founding year correction for the Voronezh university.
$ws2.Range("D7").Value = 2003

# --- Selections / active sheet ---------------------------------------------
# Update the remembered selection on the (now inactive) Университеты sheet
# before switching to and selecting a cell on the Студенты sheet, so it ends
# up as the active tab.
$null = $ws2.Range("B13").Select()
$null = $ws1.Range("D1").Select()
